$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.577
$ws.Range("B3").Value = 6.798999999999999
$ws.Range("D4").Value = -7.479000000000001
$ws.Range("B5").Value = 6.314
$ws.Range("D6").Value = -7.651999999999999
$ws.Range("C7").Value = -12.871
$ws.Range("A9").Value = -21.465
$ws.Range("C9").Value = -12.45
$ws.Range("D10").Value = -7.495
$ws.Range("B11").Value = 6.834000000000001
$ws.Range("D11").Value = -8.797000000000001
$ws.Range("B12").Value = 6.834000000000001
$ws.Range("A13").Value = -21.918
$ws.Range("A16").Value = -20.859
$ws.Range("A18").Value = -21.751
$ws.Range("A20").Value = -21.664
$ws.Range("B21").Value = 6.27
$ws.Range("C21").Value = -12.282
$ws.Range("D21").Value = -7.710000000000001
$ws.Range("D25").Value = -8.039999999999999
